$d = $word.ActiveDocument

# The document contains four "<id>p129v_N</id>" fields (N = 1..4), each
# currently split across three separate runs:
#   run 1: "<id>"     (Courier New, color 7f6000, sz 18)
#   run 2: "p129v_N"  (color 000000)
#   run 3: "</id>"    (Courier New, color 7f6000, sz 18)
#
# The edit collapses each triple into a single run containing the full
# "<id>p129v_N</id>" text. Using Find/Replace across the run boundary
# merges the found text into one run, adopting the formatting of the
# first run in the match (Courier New / 7f6000 / sz 18), which is the
# desired end state.

for ($i = 1; $i -le 4; $i++) {
    $old = "<id>p129v_$i</id>"
    $new = "<id>p129v_$i</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}
